$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.703.60'
$ws.Range('E2').Value = '  +1.51%  '

$ws.Range('D3').Value = '2.488.35'
$ws.Range('E3').Value = '  +1.72%  '

$ws.Range('E4').Value = '  -0.05%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '532.93'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +4.43%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '133.64'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +3.17%  '

$ws.Range('E7').Value = '  +0.26%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.567'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +3.12%  '

$ws.Range('D9').Value = '2.498.42'
$ws.Range('E9').Value = '  +1.42%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0990'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +3.22%  '

$ws.Range('E11').Value = '  -2.26%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '5.19'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.03%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.329'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +0.50%  '

$ws.Range('D14').Value = '2.931.03'
$ws.Range('E14').Value = '  +1.56%  '

$ws.Range('D15').Value = '58.541.73'
$ws.Range('E15').Value = '  +1.38%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '22.22'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +1.77%  '

$ws.Range('E17').Value = '  +2.03%  '

$ws.Range('D18').Value = '2.492.58'
$ws.Range('E18').Value = '  +1.17%  '

$ws.Range('E19').Value = '  +0.59%  '

$ws.Range('E20').Value = '  +2.76%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '320.03'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.68%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.22'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +5.23%  '

$ws.Range('E23').Value = '  +0.07%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '65.92'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +4.25%  '

$ws.Range('E25').Value = '  +1.59%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.993'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +0.13%  '

$ws.Range('E27').Value = '  -0.33%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '7.45'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +2.89%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '173.17'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +3.20%  '

$ws.Range('E30').Value = '  +3.77%  '

$ws.Range('E31').Value = '  +3.94%  '

$ws.Range('E32').Value = '  +2.55%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '6.25'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +1.93%  '

$ws.Range('E35').Value = '  +0.13%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '18.07'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +1.84%  '

$ws.Range('E37').Value = '  -2.77%  '

$ws.Range('E38').Value = '  +1.53%  '

$ws.Range('E39').Value = '  +3.94%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '36.26'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.80%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.808'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +6.70%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.14'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +4.41%  '

$ws.Range('E43').Value = '  +2.82%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '273.99'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.58%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '131.53'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +9.51%  '

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.591'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +1.11%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0932'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +2.11%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0510'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +4.71%  '

$ws.Range('E49').Value = '  +3.85%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '17.55'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.81%  '

$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '16.73'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.65%  '
